# FormatoPruebasRendimientoSegundoCorte.xlsx - update with new execution results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update GAP (D), Tiempo (G) and Funcion Objetivo (J/K) values for the
# "new executions" referenced in the commit message.

# Row 3 - Western Sahara
$ws.Range("D3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("J3").Value = 1.9238

# Row 4 - Djibouti
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 0.001
$ws.Range("K4").Value = 7075

# Row 5 - Qatar
$ws.Range("D5").Value = 0.003
$ws.Range("G5").Value = 0.0761
$ws.Range("J5").Value = 31.0493
$ws.Range("K5").Value = 10900

# Row 6 - Uruguay
$ws.Range("D6").Value = 0.0347
$ws.Range("G6").Value = 3.94
$ws.Range("J6").Value = 62.6178
$ws.Range("K6").Value = 95126

# Move the active selection from E8 to E7, matching the new workbook state.
$null = $ws.Range("E7").Select()
